$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5005
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 30
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 82
$ws.Range("N6").Value = -30224

$ws.Range("H8").Value = 52.75
$ws.Range("I8").Value = 57.57143
$ws.Range("J8").Value = 19
$ws.Range("K8").Value = 172.71429
$ws.Range("L8").Value = 57
$ws.Range("M8").Value = -33.71429000000001
$ws.Range("N8").Value = -335

$ws.Range("H52").Value = 1876.75
$ws.Range("I52").Value = 2002.3334
$ws.Range("K52").Value = 6007.0002
$ws.Range("M52").Value = -5847.0002

$ws.Range("H62").Value = 67553.875
$ws.Range("I62").Value = 95261.63
$ws.Range("J62").Value = 6596.8
$ws.Range("K62").Value = 95261.63
$ws.Range("L62").Value = 6596.8
$ws.Range("M62").Value = -94637.63
$ws.Range("N62").Value = -7844.8

$ws.Range("H65").Value = 67553.875
$ws.Range("I65").Value = 95261.63
$ws.Range("J65").Value = 6596.8
$ws.Range("K65").Value = 476308.15
$ws.Range("L65").Value = 32984
$ws.Range("M65").Value = -473188.15
$ws.Range("N65").Value = -39224

$ws.Range("H127").Value = 228000
$ws.Range("I127").Value = 228000
$ws.Range("K127").Value = 684000
$ws.Range("M127").Value = -679040

$ws.Range("H128").Value = 80390
$ws.Range("J128").Value = 80390
$ws.Range("L128").Value = 80390
$ws.Range("N128").Value = -90350

$ws.Range("H138").Value = 1868.9656
$ws.Range("J138").Value = 2485.4546
$ws.Range("L138").Value = 7456.3638
$ws.Range("N138").Value = -17736.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 98000
$ws.Range("J92").Value = 98000
$ws.Range("L92").Value = 98000
$ws.Range("N92").Value = -102992

$ws.Range("H132").Value = 7056.108
$ws.Range("I132").Value = 9567.68
$ws.Range("J132").Value = 1823.6666
$ws.Range("K132").Value = 28703.04
$ws.Range("L132").Value = 5470.9998
$ws.Range("M132").Value = -26173.04
$ws.Range("N132").Value = -10530.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2174.5
$ws.Range("I99").Value = 1659.4
$ws.Range("K99").Value = 1659.4
$ws.Range("M99").Value = -161.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1874.2273
$ws.Range("I31").Value = 1511.95
$ws.Range("J31").Value = 5497
$ws.Range("K31").Value = 1511.95
$ws.Range("L31").Value = 5497
$ws.Range("M31").Value = -1216.95
$ws.Range("N31").Value = -6087

$ws.Range("H34").Value = 1874.2273
$ws.Range("I34").Value = 1511.95
$ws.Range("J34").Value = 5497
$ws.Range("K34").Value = 1511.95
$ws.Range("L34").Value = 5497
$ws.Range("M34").Value = -1309.95
$ws.Range("N34").Value = -5901

$ws.Range("H122").Value = 2603.5217
$ws.Range("I122").Value = 3262.5386
$ws.Range("J122").Value = 1746.8
$ws.Range("K122").Value = 9787.6158
$ws.Range("L122").Value = 5240.4
$ws.Range("M122").Value = -7337.6158
$ws.Range("N122").Value = -10140.4

$ws.Range("H134").Value = 1197.3334
$ws.Range("I134").Value = 1197.3334
$ws.Range("K134").Value = 3592.0002
$ws.Range("M134").Value = -1057.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 128553.125
$ws.Range("I112").Value = 128553.125
$ws.Range("K112").Value = 385659.375
$ws.Range("M112").Value = -384551.375

$ws.Range("H113").Value = 420.83334
$ws.Range("J113").Value = 258.33334
$ws.Range("L113").Value = 775.0000200000001
$ws.Range("N113").Value = -5115.00002

$ws.Range("H118").Value = 93285.82000000001
$ws.Range("I118").Value = 200182.4
$ws.Range("J118").Value = 4205.3335
$ws.Range("K118").Value = 600547.2
$ws.Range("L118").Value = 12616.0005
$ws.Range("M118").Value = -599304.2
$ws.Range("N118").Value = -15102.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 21474.219
$ws.Range("I97").Value = 26811.76
$ws.Range("J97").Value = 2411.5715
$ws.Range("K97").Value = 26811.76
$ws.Range("L97").Value = 2411.5715
$ws.Range("M97").Value = -26315.76
$ws.Range("N97").Value = -3403.5715

$ws.Range("H132").Value = 3831.6538
$ws.Range("I132").Value = 3897.4211
$ws.Range("K132").Value = 11692.2633
$ws.Range("M132").Value = -9162.263300000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 888.5454999999999
$ws.Range("I22").Value = 950
$ws.Range("K22").Value = 950
$ws.Range("M22").Value = -655

$ws.Range("H27").Value = 888.5454999999999
$ws.Range("I27").Value = 950
$ws.Range("K27").Value = 950
$ws.Range("M27").Value = -843

$ws.Range("H46").Value = 2249.625
$ws.Range("I46").Value = 2399
$ws.Range("J46").Value = 2000.6666
$ws.Range("K46").Value = 2399
$ws.Range("L46").Value = 2000.6666
$ws.Range("M46").Value = -2211
$ws.Range("N46").Value = -2376.6666

$ws.Range("H93").Value = 1550.4736
$ws.Range("I93").Value = 1469.9445
$ws.Range("K93").Value = 1469.9445
$ws.Range("M93").Value = -221.9445000000001

$ws.Range("H100").Value = 6988.5557
$ws.Range("I100").Value = 8599.333000000001
$ws.Range("J100").Value = 3767
$ws.Range("K100").Value = 8599.333000000001
$ws.Range("L100").Value = 3767
$ws.Range("M100").Value = -8058.333000000001
$ws.Range("N100").Value = -4849

$ws.Range("H122").Value = 3893.0715
$ws.Range("I122").Value = 3625.4167
$ws.Range("K122").Value = 10876.2501
$ws.Range("M122").Value = -8426.250100000001

$ws.Range("H136").Value = 3035.1875
$ws.Range("I136").Value = 2051.6365
$ws.Range("J136").Value = 5199
$ws.Range("K136").Value = 6154.9095
$ws.Range("L136").Value = 15597
$ws.Range("M136").Value = -3604.9095
$ws.Range("N136").Value = -20697

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws.Range("H61").Value = 26164.666
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null

$ws.Range("H62").Value = 6814.625
$ws.Range("J62").Value = 5972.3335
$ws.Range("L62").Value = 5972.3335
$ws.Range("N62").Value = -7220.3335

$ws.Range("H65").Value = 6814.625
$ws.Range("J65").Value = 5972.3335
$ws.Range("L65").Value = 29861.6675
$ws.Range("N65").Value = -36101.6675

$ws.Range("H96").Value = 1633.1666
$ws.Range("I96").Value = 1759.8
$ws.Range("K96").Value = 1759.8
$ws.Range("M96").Value = -386.8
